$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column AR: header "03-ago" continuing the date series in row 1,
# plus the corresponding data values for rows 2-18.
$ws.Range("AR1").Value = "03-ago"

$values = @{
    2  = 0
    3  = 15.689327096616635
    4  = 19.966936274338202
    5  = 18.822531355581052
    6  = 0
    7  = 15.114920988675921
    8  = 8.3982846636579307
    9  = 12.541465128871852
    10 = 11.648992059434557
    11 = 15.023644843967348
    12 = 0
    13 = 6.6675980200815719
    14 = 0
    15 = 0
    16 = 11.506558368047143
    17 = 0
    18 = 0
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 44).Value = $values[$row]
}

# Replicate the author's final cursor/selection position.
[void]$ws.Range("AT7").Select()
